# msz - 3./4. smoke test + inheritance page and 2. dialog
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New layout of the "pages" table (rows 2-25). Rows 2-13 are overwritten in
# place; rows 14-25 are newly appended below the previous last row (13).

$rows = @(
    @("dlgAutomobileInsurance_pagVehicleData", $null),
    @("dlgTruckInsurance_pagVehicleData", $null),
    @("dlgMotorcycleInsurance_pagVehicleData", $null),
    @("dlgCamperInsurance_pagVehicleData", $null),
    @("dlgAutomobileInsurance_pagInsurantData", "dlgBaseInsurance_pagInsurantData"),
    @("dlgTruckInsurance_pagInsurantData", "dlgBaseInsurance_pagInsurantData"),
    @("dlgMotorcycleInsurance_pagInsurantData", "dlgBaseInsurance_pagInsurantData"),
    @("dlgCamperInsurance_pagInsurantData", "dlgBaseInsurance_pagInsurantData"),
    @("dlgAutomobileInsurance_pagProductData", $null),
    @("dlgTruckInsurance_pagProductData", $null),
    @("dlgMotorcycleInsurance_pagProductData", $null),
    @("dlgCamperInsurance_pagProductData", $null),
    @("dlgAutomobileInsurance_pagPriceOption", "dlgBaseInsurance_pagPriceOption"),
    @("dlgTruckInsurance_pagPriceOption", "dlgBaseInsurance_pagPriceOption"),
    @("dlgMotorcycleInsurance_pagPriceOption", "dlgBaseInsurance_pagPriceOption"),
    @("dlgCamperInsurance_pagPriceOption", "dlgBaseInsurance_pagPriceOption"),
    @("dlgAutomobileInsurance_pagSendQuote", "dlgBaseInsurance_pagSendQuote"),
    @("dlgTruckInsurance_pagSendQuote", "dlgBaseInsurance_pagSendQuote"),
    @("dlgMotorcycleInsurance_pagSendQuote", "dlgBaseInsurance_pagSendQuote"),
    @("dlgCamperInsurance_pagSendQuote", "dlgBaseInsurance_pagSendQuote"),
    @("dlgMain_pagTrainings", $null),
    @("dlgProfil_pagAbwesenheiten", $null),
    @("dlgProfil_pagBenachrichtigungen", $null),
    @("dlgProfil_pagProfil", $null)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1]) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        $ws.Cells.Item($r, 2).Value = $null
    }
    $r = $r + 1
}

$ws.Range("B28").Select() | Out-Null
